# Applies the "Updated abstract and introduction" edit to Introduction.docx
$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Replace the body text of paragraph 2 (the big "Our team..." one)
# -----------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$rng2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$rng2.Text = "Our team has been given the task to develop a booking system for the company VIA Bus. An interview was conducted with the owners of the company to decide what kind of system they wanted.  It was concluded that the system must allow the employees to handle reservations for seats for trips and journeys to different premade destinations or a whole bus with personal preferences (food, party guide, additional stops) to desired destination. The system contains information about fixed tours (trips and journeys), and non-fixed tours(bus-and-chauffeurs) as well as chauffeurs, customers, and passenger’s data. The company told our team that only the employees will be given access to the system and the program won’t be accessible to the public."

# Re-insert the _GoBack bookmark: it now sits between "concluded tha" and "t the system"
$p2 = $d.Paragraphs(2)
$bmRng = $p2.Range.Duplicate
$bmRng.Find.ClearFormatting()
$bmRng.Find.Execute("concluded tha") | Out-Null
$bmRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

# -----------------------------------------------------------------
# 2) Replace the body text of paragraph 3 (UI / Java paragraph)
# -----------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$rng3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$rng3.Text = "The system has a simple user interface with interactive elements and is developed using Java programming language. To meet the owner’s requirements, we’ve implemented a way to create/edit/remove different types of tours, chauffeurs, customers, and reservations. Passengers can be added and removed, but they cannot be edited. We’ve also added a way for the company to keep track of the number of busses they have, so that every time a bus is used, the number goes down and when the bus is no longer being used, the number goes up."

# -----------------------------------------------------------------
# 3) Replace the body text of paragraph 4 (purpose paragraph)
# -----------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$rng4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$rng4.Text = "The purpose of this system is to make it easier and more efficient for VIA Bus to manage all the reservations and keep track of all the data that the company works with. In the following pages of the report, we will go into greater detail about the things we’ve done to create a simple, yet efficient system - creating a proper system analysis, containing the requirements, detailed activity diagrams and use cases, also a detailed GUI breakdown as well as an UML breakdown was made so that we can give the thoughts that we had while making the system. Implementation was also added to provide some details on how we approached different tasks code-wise. We also went through several different test scenarios so that we can establish if the system is properly working. "

# -----------------------------------------------------------------
# 4) Formatting: font size + paragraph spacing for the first four
#    paragraphs (title + the three body paragraphs). The trailing
#    blank paragraph keeps its original 34/34 size and spacing.
# -----------------------------------------------------------------
for ($i = 1; $i -le 4; $i++) {
    $p = $d.Paragraphs($i)
    $p.Range.Font.Size = 12
    $p.Range.Font.SizeBi = 12
    $p.Range.ParagraphFormat.SpaceAfter = 0
}

Write-Host "done"
